# Adds a new "2022-Q3" sheet (as the 2nd tab) with fund-holding detail data,
# and updates the "总计" (summary) sheet on tab 1 with a new 2022-Q3 row,
# pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Forces a numeric-looking string (e.g. "005613", "0.1432") to be stored
    # as text rather than being auto-coerced into a number by the COM layer.
    param($ws, $row, $col, $val)
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new 2022-Q3 row at the top of the data
#    (row 2), shifting 2022-Q2 / 2022-Q1 / 2021-Q3 / 2021-Q2 down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Extend the styled index column (A) down to the new last row (row 6) by
# copying the formatting already used by the rows above it.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

# Write bottom-up so we never overwrite a row before reading/using its data.
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q2"
$summary.Cells.Item(6, 3).Value = 7
$summary.Cells.Item(6, 4).Value = 0.53

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2021-Q3"
$summary.Cells.Item(5, 3).Value = 7
$summary.Cells.Item(5, 4).Value = 0.52

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 8
$summary.Cells.Item(4, 4).Value = 0.75

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 9
$summary.Cells.Item(3, 4).Value = 0.74

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 11
$summary.Cells.Item(2, 4).Value = 0.63

# ---------------------------------------------------------------------------
# 2) Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" sheet, so tab order becomes:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q3, 2021-Q2
# ---------------------------------------------------------------------------
$q2SheetForInsert = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2SheetForInsert)
$q3.Name = "2022-Q3"

# NOTE: the old $q2SheetForInsert reference goes stale the moment a new sheet
# is inserted (its Font/format reads start coming back wrong) - always grab
# a fresh reference by name right before reading from it.

# Match page margins used by the other quarterly sheets.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Match the outline summary settings used by the other sheets.
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

# Copy the header-row formatting (bold + border, style used by B1:H1 on the
# other quarterly sheets) from the existing "2022-Q2" sheet.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A index-cell formatting (bold + border) too.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("A2").Copy()
$q3.Range("A2:A12").PasteSpecial(-4122)

# Header row
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Data rows: index, code, name, size, position, ratio, marketValue, rank
$rows = @(
    @(0, "005613", "上投摩根富时发达市场REITs指数（QDII）人民币份额", "3.41", "92.80", "4.20", "0.1432", 3),
    @(1, "005614", "上投摩根富时发达市场REITs指数（QDII）美钞", "3.41", "92.80", "4.20", "0.1432", 3),
    @(2, "005615", "上投摩根富时发达市场REITs指数（QDII）美汇", "3.41", "92.80", "4.20", "0.1432", 3),
    @(3, "000179", "广发美国房地产指数（QDII）人民币A", "1.82", "92.37", "3.44", "0.0626", 4),
    @(4, "000180", "广发美国房地产指数（QDII）美元A", "1.82", "92.37", "3.44", "0.0626", 4),
    @(5, "160140", "南方道琼斯美国精选REIT指数（QDII-LOF）A", "0.78", "91.13", "3.86", "0.0301", 4),
    @(6, "070031", "嘉实全球房地产（QDII）", "0.38", "94.39", "4.47", "0.0170", 2),
    @(7, "160141", "南方道琼斯美国精选REIT指数（QDII-LOF）C", "0.42", "91.13", "3.86", "0.0162", 4),
    @(8, "320017", "诺安全球收益不动产（QDII）", "0.23", "73.76", "5.43", "0.0125", 4),
    @(9, "016278", "广发美国房地产指数（QDII）人民币C", "0.01", "92.37", "3.44", "0.0003", 4),
    @(10, "016279", "广发美国房地产指数（QDII）美元C", "0.01", "92.37", "3.44", "0.0003", 4)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $data[0]
    Set-TextValue $q3 $r 2 $data[1]
    $q3.Cells.Item($r, 3).Value = $data[2]
    Set-TextValue $q3 $r 4 $data[3]
    Set-TextValue $q3 $r 5 $data[4]
    Set-TextValue $q3 $r 6 $data[5]
    Set-TextValue $q3 $r 7 $data[6]
    $q3.Cells.Item($r, 8).Value = $data[7]
}
